# Update "想去人数" (interest/attendance count) figures in the 展览 (sheet1)
# and 全部类型 (sheet4) worksheets to match newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - column F holds the counts for rows 3,4,7,8,11,13,14,15,22
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 72
$wsExhibit.Range("F4").Value = 1508
$wsExhibit.Range("F7").Value = 11072
$wsExhibit.Range("F8").Value = 11072
$wsExhibit.Range("F11").Value = 321
$wsExhibit.Range("F13").Value = 755
$wsExhibit.Range("F14").Value = 12229
$wsExhibit.Range("F15").Value = 12773
$wsExhibit.Range("F22").Value = 40

# Sheet "全部类型" - same events, shifted down by one row (rows 4,5,8,9,12,14,15,16,23)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 72
$wsAll.Range("F5").Value = 1508
$wsAll.Range("F8").Value = 11072
$wsAll.Range("F9").Value = 11072
$wsAll.Range("F12").Value = 321
$wsAll.Range("F14").Value = 755
$wsAll.Range("F15").Value = 12229
$wsAll.Range("F16").Value = 12773
$wsAll.Range("F23").Value = 40
